# Weekly driver report update for 2025-04-29
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# (ColumnWidth setter adds a constant 5/6-character padding on this engine,
#  so back that out to land on the exact target widths stored in the XML.)
$ws.Columns.Item(2).ColumnWidth  = 14 - (5/6)   # B: 15 -> 14
$ws.Columns.Item(5).ColumnWidth  = 14 - (5/6)   # E: 16 -> 14
$ws.Columns.Item(6).ColumnWidth  = 11 - (5/6)   # F: 2  -> 11
$ws.Columns.Item(7).ColumnWidth  = 48 - (5/6)   # G: 2  -> 48
$ws.Columns.Item(8).ColumnWidth  = 15 - (5/6)   # H: 2  -> 15
$ws.Columns.Item(9).ColumnWidth  = 30 - (5/6)   # I: 2  -> 30
$ws.Columns.Item(10).ColumnWidth = 16 - (5/6)   # J: 2  -> 16

# --- Updated "Bad Drivers" figures (row 3 / row 4 totals) ---
$ws.Range("C3").Value = 130
$ws.Range("D3").Value = 98.59999999999999
$ws.Range("C4").Value = 130

# --- Rebuild the "Good Drivers" header row (row 11) with new columns,
#     stripped of the old bold/border header style ---
$ws.Range("A11:J11").ClearFormats()

$ws.Range("A11").Value = "adapter-driver"
$ws.Range("B11").Value = "good sum"
$ws.Range("C11").Value = "critical sum"
$ws.Range("D11").Value = "warning sum"
$ws.Range("E11").Value = "client count"
$ws.Range("F11").Value = "total sum"
$ws.Range("G11").Value = "adapter"
$ws.Range("H11").Value = "driver"
$ws.Range("I11").Value = "good roaming calculation (%)"
$ws.Range("J11").Value = "driver vintage"

# --- New data row 12 (replaces previously-empty rows 12-16) ---
$ws.Range("A12").Value = "Realtek RTL8852AE WiFi 6 802.11ax PCIe Adapter - 6001.10.356.1"
$ws.Range("B12").Value = 1071383
$ws.Range("C12").Value = 4419
$ws.Range("D12").Value = 180
$ws.Range("E12").Value = 1644
$ws.Range("F12").Value = 1075982
$ws.Range("G12").Value = "realtek rtl8852ae wifi 6 802.11ax pcie adapter"
$ws.Range("H12").Value = "6001.10.356.1"
$ws.Range("I12").Value = 99.59999999999999

# Force the vintage date to be stored as literal text, not an auto-parsed date
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = "2024-05-12"
$ws.Range("J12").ClearFormats()
